$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark (it currently sits in the
#    "Thanks for this suggestion!" paragraph) so we can re-create it
#    in its new location without a name collision.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Append ". See Screen." right after "...Learn Git" (end of the
#    paragraph that currently reads "Completed Code Academy Course
#    Learn Git"), then place an empty "_GoBack" bookmark right after
#    the newly inserted text.
# ------------------------------------------------------------------
$gitPara = $d.Paragraphs.Item(9)
$gitRange = $d.Range($gitPara.Range.Start, $gitPara.Range.End - 1)
$gitRange.InsertAfter(". See Screen.")

# Re-fetch the (now longer) paragraph end, insert a throw-away marker
# character, wrap a bookmark around it, then delete the marker again.
# Deleting the wrapped character leaves the bookmark collapsed at the
# correct spot (Word collapses a bookmark when its content is deleted).
$paraEnd = $d.Paragraphs.Item(9).Range.End - 1
$markerRange = $d.Range($paraEnd, $paraEnd)
$markerRange.InsertAfter("X")
$bookmarkRange = $d.Range($paraEnd, $paraEnd + 1)
$bookmarkRange.Bookmarks.Add("_GoBack")
$deleteRange = $d.Range($paraEnd, $paraEnd + 1)
$deleteRange.Delete()

# ------------------------------------------------------------------
# 3) Merge the three separate runs of the "Overall, ..." paragraph
#    into a single run (same combined text, no run-splitting).
# ------------------------------------------------------------------
$overallPara = $d.Paragraphs.Item(15)
$overallRange = $d.Range($overallPara.Range.Start, $overallPara.Range.End - 1)
$overallText = $overallRange.Text
$overallRange.Delete()
$insertionPoint = $d.Range($overallPara.Range.Start, $overallPara.Range.Start)
$insertionPoint.InsertAfter($overallText)
